$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from an existing header cell (A1) onto the
# three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the win/loss/tie record for each data row (2 through 44)
for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 30).Value = 79
    $ws.Cells.Item($row, 31).Value = 83
    $ws.Cells.Item($row, 32).Value = 0
}
